# CompanyMaster.xlsx - add new login-credential test rows and update the
# "Login" button locator test step on the LoginTest sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginCredentials")
$ws2 = $wb.Worksheets.Item("LoginTest")

# --- LoginTest: the "Login" step now locates the button by className instead
#     of its old (wrong) id value ---
$ws2.Range("D4").Value = "className"
$ws2.Range("E4").Value = "lfr-btn-label"

# --- LoginCredentials: append two more sample rows under the existing pair ---
$ws1.Range("A3").Value = "john.smith"
$ws1.Range("B3").Value = "test1"
$ws1.Range("A4").Value = "john.smith1"
$ws1.Range("B4").Value = "askd"

# --- restore the window size/position recorded the last time the workbook
#     was saved, and leave the selection where the author left it ---
$win = $excel.ActiveWindow
$win.Left   = 2055
$win.Top    = 4155
$win.Width  = 14400
$win.Height = 10755

$ws1.Activate()
$ws1.Range("B4").Select()
